# Apply weekly refresh of Fruta/Hortaliza data: rows 2-8 get their
# variable columns (D, L, M, N, O, P, Q, R, S, T) permuted according to a
# single 7-cycle: 2->7->3->5->8->6->4->2
# (i.e. the data that used to live in row 7 now lives in row 2, the data
# that used to live in row 5 now lives in row 3, etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values move together as a "record" for this weekly swap.
$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# Capture the current ("before") values for every affected row/column
# so we can rewrite them in their new positions without clobbering data
# we still need to read later in the loop.
$rows = @(2, 3, 4, 5, 6, 7, 8)
$before = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$r").Value2
    }
    $before[$r] = $rowData
}

# Destination row -> source row (the row whose old data it should receive).
$mapping = @{
    2 = 7
    3 = 5
    4 = 2
    5 = 8
    6 = 4
    7 = 3
    8 = 6
}

foreach ($destRow in $rows) {
    $srcRow = $mapping[$destRow]
    $srcData = $before[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $srcData[$col]
    }
}
